$d = $word.ActiveDocument
$ps = $d.PageSetup
$ps.PageWidth = 595.3
$ps.PageHeight = 841.9
$ps.TopMargin = 70.85
$ps.RightMargin = 70.85
$ps.BottomMargin = 70.85
$ps.LeftMargin = 70.85
$ps.HeaderDistance = 35.4
$ps.FooterDistance = 35.4
$ps.Gutter = 0
